$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.575.62'
$ws.Range('E2').Value = '  +0.61%  '

$ws.Range('D3').Value = '1.477.94'
$ws.Range('E3').Value = '  +0.82%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.012'
$ws.Range('E4').Value = '  +0.18%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9607'
$ws.Range('E5').Value = '  +4.53%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '280.61'
$ws.Range('E6').Value = '  +0.02%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3663'
$ws.Range('E7').Value = '  -1.36%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3075'
$ws.Range('E8').Value = '  -3.54%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '40.02'
$ws.Range('E9').Value = '  -1.25%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.062'
$ws.Range('E10').Value = '  +0.80%  '

$ws.Range('E11').Value = '  +0.29%  '

$ws.Range('E12').Value = '  +0.04%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.521'
$ws.Range('E13').Value = '  -0.86%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.06'
$ws.Range('E14').Value = '  -0.33%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.215'
$ws.Range('E15').Value = '  -0.10%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9617'
$ws.Range('E16').Value = '  +3.67%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001034'
$ws.Range('E17').Value = '  -0.09%  '

$ws.Range('D18').Value = '1.476.93'
$ws.Range('E18').Value = '  +0.01%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.05972'
$ws.Range('E19').Value = '  +4.30%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.96'
$ws.Range('E20').Value = '  -2.23%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.505'
$ws.Range('E21').Value = '  -3.32%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.45'
$ws.Range('E22').Value = '  -1.69%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.05'
$ws.Range('E23').Value = '  -1.23%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.256'
$ws.Range('E24').Value = '  -1.56%  '

$ws.Range('D25').Value = '20.635.34'
$ws.Range('E25').Value = '  +0.11%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '143.41'
$ws.Range('E26').Value = '  +3.84%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.111'
$ws.Range('E27').Value = '  -8.20%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.28'
$ws.Range('E28').Value = '  -1.38%  '

$ws.Range('D29').Value = '1.639.04'
$ws.Range('E29').Value = '  +0.20%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '113.85'
$ws.Range('E30').Value = '  +0.30%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.982'
$ws.Range('E31').Value = '  +0.33%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.029'
$ws.Range('E32').Value = '  -4.79%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8164'
$ws.Range('E33').Value = '  -3.61%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07975'
$ws.Range('E34').Value = '  +2.24%  '

$ws.Range('E35').Value = '  -1.72%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.222'
$ws.Range('E36').Value = '  +6.41%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05797'
$ws.Range('E37').Value = '  -4.83%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.729'
$ws.Range('E38').Value = '  -3.08%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02050'
$ws.Range('E39').Value = '  -0.69%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9613'
$ws.Range('E40').Value = '  +2.12%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.39'
$ws.Range('E41').Value = '  -2.49%  '

$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1878'
$ws.Range('E42').Value = '  -0.82%  '

$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.474'
$ws.Range('E43').Value = '  +1.39%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5313'
$ws.Range('E44').Value = '  -1.44%  '

$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.546'
$ws.Range('E45').Value = '  -1.22%  '

$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.23'
$ws.Range('E46').Value = '  -1.01%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '118.36'
$ws.Range('E47').Value = '  -4.34%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5202'
$ws.Range('E48').Value = '  -1.94%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.821'
$ws.Range('E49').Value = '  -0.31%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06495'
$ws.Range('E50').Value = '  +0.79%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9912'
$ws.Range('E51').Value = '  -0.31%  '
